$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM fixes: remove duplicate-part confusion by updating two rows ---

# Row 21 (AMS1117-3.3): manufacturer + LCSC part number corrected
# Set the LCSC part number first, then the manufacturer, so new shared
# strings are appended to the string table in that order.
$ws.Range("F21").Value = "C347222"
$ws.Range("D21").Value = "UMW(Youtai Semiconductor Co., Ltd.)"

# Row 25: manufacture part number + LCSC part number corrected
$ws.Range("F25").Value = "C94674"
$ws.Range("C25").Value = "X1A000141000100"

# Column D widened slightly to fit the new, longer manufacturer name
$ws.Columns("D").ColumnWidth = 33.17

# Active selection left on C25 (last cell touched)
$ws.Range("C25").Select()
